$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "tropical" subject to "tropical medicine" everywhere it appears
# on the sheet (affects every row using that shared string, e.g. the Subject
# column for groups B2A-B2E).
$ws.Cells.Replace("tropical", "tropical medicine", 1, 1, $false, $false, $false)
